$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was row 6's data)
$ws.Range("A2").Value = 112438849
$ws.Range("B2").Value = 95704
$ws.Range("E2").Value = 221946
$ws.Range("F2").Value = "Mattlummer"
$ws.Range("G2").Value = "Lycopodium clavatum"
$ws.Range("Q2").Value = 502192
$ws.Range("R2").Value = 6543228

# Row 3 (was row 5's data)
$ws.Range("A3").Value = 112438845
$ws.Range("B3").Value = 95704
$ws.Range("E3").Value = 221946
$ws.Range("F3").Value = "Mattlummer"
$ws.Range("G3").Value = "Lycopodium clavatum"
$ws.Range("Q3").Value = 502317
$ws.Range("R3").Value = 6543245

# Row 4 (was row 3's data)
$ws.Range("A4").Value = 112438846
$ws.Range("B4").Value = 95704
$ws.Range("E4").Value = 221946
$ws.Range("F4").Value = "Mattlummer"
$ws.Range("G4").Value = "Lycopodium clavatum"
$ws.Range("Q4").Value = 502269
$ws.Range("R4").Value = 6543231

# Row 5 (was row 2's data)
$ws.Range("A5").Value = 112438847
$ws.Range("B5").Value = 95707
$ws.Range("E5").Value = 221941
$ws.Range("F5").Value = "Plattlummer"
$ws.Range("G5").Value = "Lycopodium complanatum"
$ws.Range("Q5").Value = 502260
$ws.Range("R5").Value = 6543183

# Row 6 (was row 4's data)
$ws.Range("A6").Value = 112438848
$ws.Range("B6").Value = 95707
$ws.Range("E6").Value = 221941
$ws.Range("F6").Value = "Plattlummer"
$ws.Range("G6").Value = "Lycopodium complanatum"
$ws.Range("Q6").Value = 502199
$ws.Range("R6").Value = 6543178
